# Update positions of balises in TrackDescription (column A holds the
# balise position in meters). Each value below gains a trailing precision
# digit, matching the source-controlled change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 241
$ws.Range("A3").Value = 1181
$ws.Range("A12").Value = 2801
$ws.Range("A13").Value = 3612
$ws.Range("A14").Value = 3666
$ws.Range("A15").Value = 3812
$ws.Range("A16").Value = 3868
$ws.Range("A17").Value = 4043
$ws.Range("A18").Value = 4215
$ws.Range("A19").Value = 4270
$ws.Range("A20").Value = 4695
$ws.Range("A21").Value = 4749
$ws.Range("A22").Value = 4990
$ws.Range("A23").Value = 5044
$ws.Range("A24").Value = 5213
$ws.Range("A25").Value = 5268
$ws.Range("A26").Value = 5420
$ws.Range("A27").Value = 5564

# Scroll the view down and move the selection, matching the saved
# sheetView state (topLeftCell A13->A22, selection A25->A28).
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 22
$ws.Range("A28").Select()
